$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 8
$ws_ALC.Range("H8").Value = 465.54544
$ws_ALC.Range("I8").Value = 235.55556
$ws_ALC.Range("J8").Value = 1500.5
$ws_ALC.Range("K8").Value = 706.66668
$ws_ALC.Range("L8").Value = 4501.5
$ws_ALC.Range("M8").Value = -567.66668
$ws_ALC.Range("N8").Value = -4779.5

# ALC row 9
$ws_ALC.Range("H9").Value = 690.125
$ws_ALC.Range("I9").Value = 770.8570999999999
$ws_ALC.Range("K9").Value = 770.8570999999999
$ws_ALC.Range("M9").Value = -601.8570999999999

# ALC row 12
$ws_ALC.Range("H12").Value = 650.1
$ws_ALC.Range("I12").Value = 611.2222
$ws_ALC.Range("K12").Value = 611.2222
$ws_ALC.Range("M12").Value = -441.2222

# ALC row 15
$ws_ALC.Range("H15").Value = 1878.6129
$ws_ALC.Range("I15").Value = 1878.6129
$ws_ALC.Range("K15").Value = 5635.8387
$ws_ALC.Range("M15").Value = -5466.8387

# ALC row 40
$ws_ALC.Range("H40").Value = 1987.1666
$ws_ALC.Range("J40").Value = 2149.5
$ws_ALC.Range("L40").Value = 2149.5
$ws_ALC.Range("N40").Value = -2499.5

# ALC row 43
$ws_ALC.Range("H43").Value = 7119
$ws_ALC.Range("I43").Value = 5501
$ws_ALC.Range("J43").Value = 7658.3335
$ws_ALC.Range("K43").Value = 5501
$ws_ALC.Range("L43").Value = 7658.3335
$ws_ALC.Range("M43").Value = -5432
$ws_ALC.Range("N43").Value = -7796.3335

# ALC row 88
$ws_ALC.Range("H88").Value = 5208.8
$ws_ALC.Range("J88").Value = 5536
$ws_ALC.Range("L88").Value = 5536
$ws_ALC.Range("N88").Value = -6348

# ALC row 91
$ws_ALC.Range("H91").Value = 5208.8
$ws_ALC.Range("J91").Value = 5536
$ws_ALC.Range("L91").Value = 5536
$ws_ALC.Range("N91").Value = -8344

# ALC row 92
$ws_ALC.Range("H92").Value = 873.1
$ws_ALC.Range("I92").Value = 725.7778
$ws_ALC.Range("K92").Value = 725.7778
$ws_ALC.Range("M92").Value = 522.2222

# ALC row 106
$ws_ALC.Range("H106").Value = 2939.25
$ws_ALC.Range("I106").Value = 2922.3333
$ws_ALC.Range("K106").Value = 2922.3333
$ws_ALC.Range("M106").Value = -2291.3333

# ALC row 107
$ws_ALC.Range("H107").Value = 1269.5714
$ws_ALC.Range("I107").Value = 1266.8572
$ws_ALC.Range("J107").Value = 1275
$ws_ALC.Range("K107").Value = 1266.8572
$ws_ALC.Range("L107").Value = 1275
$ws_ALC.Range("M107").Value = 653.1428000000001
$ws_ALC.Range("N107").Value = -5115

# ALC row 116
$ws_ALC.Range("H116").Value = 4468.294
$ws_ALC.Range("I116").Value = 4024.182
$ws_ALC.Range("J116").Value = 5282.5
$ws_ALC.Range("K116").Value = 4024.182
$ws_ALC.Range("L116").Value = 5282.5
$ws_ALC.Range("M116").Value = -582.1819999999998
$ws_ALC.Range("N116").Value = -12166.5

# ALC row 137
$ws_ALC.Range("H137").Value = 3513.9583
$ws_ALC.Range("J137").Value = 4264.933
$ws_ALC.Range("L137").Value = 12794.799
$ws_ALC.Range("N137").Value = -17894.799

# ALC row 138
$ws_ALC.Range("H138").Value = 3196.6875
$ws_ALC.Range("J138").Value = 3628.9167
$ws_ALC.Range("L138").Value = 10886.7501
$ws_ALC.Range("N138").Value = -21166.7501

# ARM row 32
$ws_ARM.Range("H32").Value = 8214.029
$ws_ARM.Range("I32").Value = 8214.029
$ws_ARM.Range("K32").Value = 8214.029
$ws_ARM.Range("M32").Value = -7927.029

# ARM row 74
$ws_ARM.Range("H74").Value = 1791.25
$ws_ARM.Range("I74").Value = 1388.5834
$ws_ARM.Range("K74").Value = 1388.5834
$ws_ARM.Range("M74").Value = -514.5834

# ARM row 77
$ws_ARM.Range("H77").Value = 1791.25
$ws_ARM.Range("I77").Value = 1388.5834
$ws_ARM.Range("K77").Value = 6942.916999999999
$ws_ARM.Range("M77").Value = -2574.916999999999

# ARM row 122
$ws_ARM.Range("H122").Value = 1863
$ws_ARM.Range("I122").Value = 1863
$ws_ARM.Range("K122").Value = 5589
$ws_ARM.Range("M122").Value = -3139

# CRP row 17
$ws_CRP.Range("H17").Value = 9975
$ws_CRP.Range("I17").Value = 10000
$ws_CRP.Range("K17").Value = 10000
$ws_CRP.Range("M17").Value = -9826

# CRP row 41
$ws_CRP.Range("H41").Value = 14053.375
$ws_CRP.Range("J41").Value = 34890
$ws_CRP.Range("L41").Value = 34890
$ws_CRP.Range("N41").Value = -35746

# CRP row 50
$ws_CRP.Range("H50").Value = 15000
$ws_CRP.Range("I50").Value = 15000
$ws_CRP.Range("K50").Value = 15000
$ws_CRP.Range("M50").Value = -14375

# CRP row 93
$ws_CRP.Range("H93").Value = 5203.5
$ws_CRP.Range("I93").Value = 5203.5
$ws_CRP.Range("K93").Value = 5203.5
$ws_CRP.Range("M93").Value = -3331.5

# CRP row 105
$ws_CRP.Range("H105").Value = 834.25
$ws_CRP.Range("I105").Value = 769
$ws_CRP.Range("J105").Value = 1030
$ws_CRP.Range("K105").Value = 769
$ws_CRP.Range("L105").Value = 1030
$ws_CRP.Range("M105").Value = 978
$ws_CRP.Range("N105").Value = -4524

# CRP row 122
$ws_CRP.Range("H122").Value = 4329.778
$ws_CRP.Range("I122").Value = 4596.4
$ws_CRP.Range("J122").Value = 3996.5
$ws_CRP.Range("K122").Value = 13789.2
$ws_CRP.Range("L122").Value = 11989.5
$ws_CRP.Range("M122").Value = -11339.2
$ws_CRP.Range("N122").Value = -16889.5

# CRP row 134
$ws_CRP.Range("H134").Value = 2071.8572
$ws_CRP.Range("I134").Value = 2174.5833
$ws_CRP.Range("K134").Value = 6523.749899999999
$ws_CRP.Range("M134").Value = -3988.749899999999

# CUL row 5
$ws_CUL.Range("H5").Value = 3738.4443
$ws_CUL.Range("I5").Value = 2003.909
$ws_CUL.Range("J5").Value = 4930.9375
$ws_CUL.Range("K5").Value = 6011.727000000001
$ws_CUL.Range("L5").Value = 14792.8125
$ws_CUL.Range("M5").Value = -5899.727000000001
$ws_CUL.Range("N5").Value = -15016.8125

# CUL row 8
$ws_CUL.Range("H8").Value = 4027.3333
$ws_CUL.Range("I8").Value = 4027.3333
$ws_CUL.Range("K8").Value = 12081.9999
$ws_CUL.Range("M8").Value = -11942.9999

# CUL row 68
$ws_CUL.Range("H68").Value = 4512.5
$ws_CUL.Range("J68").Value = 4512.5
$ws_CUL.Range("L68").Value = 13537.5
$ws_CUL.Range("N68").Value = -15159.5

# CUL row 71
$ws_CUL.Range("H71").Value = 4512.5
$ws_CUL.Range("J71").Value = 4512.5
$ws_CUL.Range("L71").Value = 40612.5
$ws_CUL.Range("N71").Value = -48724.5

# CUL row 134
$ws_CUL.Range("H134").Value = 11001.9
$ws_CUL.Range("I134").Value = 1109.5555
$ws_CUL.Range("K134").Value = 3328.6665
$ws_CUL.Range("M134").Value = 1741.3335

# CUL row 135
$ws_CUL.Range("H135").Value = 3738.4443
$ws_CUL.Range("I135").Value = 2003.909
$ws_CUL.Range("J135").Value = 4930.9375
$ws_CUL.Range("K135").Value = 18035.181
$ws_CUL.Range("L135").Value = 44378.4375
$ws_CUL.Range("M135").Value = -15500.181
$ws_CUL.Range("N135").Value = -49448.4375

# CUL row 138
$ws_CUL.Range("H138").Value = 14410
$ws_CUL.Range("I138").Value = 14410
$ws_CUL.Range("K138").Value = 43230
$ws_CUL.Range("M138").Value = -38090

# CUL row 139
$ws_CUL.Range("H139").Value = 2693.1667
$ws_CUL.Range("I139").Value = 2693.1667
$ws_CUL.Range("K139").Value = 8079.500100000001
$ws_CUL.Range("M139").Value = -2939.500100000001

# GSM row 2
$ws_GSM.Range("H2").Value = 1000.5
$ws_GSM.Range("I2").Value = 400.6
$ws_GSM.Range("K2").Value = 400.6
$ws_GSM.Range("M2").Value = -287.6

# GSM row 122
$ws_GSM.Range("H122").Value = 1924
$ws_GSM.Range("I122").Value = 1924
$ws_GSM.Range("J122").Value = 0
$ws_GSM.Range("K122").Value = 5772
$ws_GSM.Range("L122").Value = 0
$ws_GSM.Range("M122").Value = -3322
$ws_GSM.Range("N122").ClearContents()

# LTW row 7
$ws_LTW.Range("H7").Value = 4561.5713
$ws_LTW.Range("I7").Value = 4386.4
$ws_LTW.Range("K7").Value = 4386.4
$ws_LTW.Range("M7").Value = -4274.4

# LTW row 40
$ws_LTW.Range("H40").Value = 9999.333000000001
$ws_LTW.Range("I40").Value = 1999.2
$ws_LTW.Range("K40").Value = 1999.2
$ws_LTW.Range("M40").Value = -1863.2

# LTW row 82
$ws_LTW.Range("H82").Value = 3246.7693
$ws_LTW.Range("J82").Value = 4811.3335
$ws_LTW.Range("L82").Value = 4811.3335
$ws_LTW.Range("N82").Value = -5533.3335

# LTW row 85
$ws_LTW.Range("H85").Value = 3246.7693
$ws_LTW.Range("J85").Value = 4811.3335
$ws_LTW.Range("L85").Value = 4811.3335
$ws_LTW.Range("N85").Value = -7307.3335

# LTW row 100
$ws_LTW.Range("H100").Value = 7861.8
$ws_LTW.Range("I100").Value = 8101.6665
$ws_LTW.Range("K100").Value = 8101.6665
$ws_LTW.Range("M100").Value = -7560.6665

# LTW row 122
$ws_LTW.Range("H122").Value = 0
$ws_LTW.Range("I122").Value = 0
$ws_LTW.Range("J122").Value = 0
$ws_LTW.Range("K122").Value = 0
$ws_LTW.Range("L122").Value = 0
$ws_LTW.Range("M122").ClearContents()
$ws_LTW.Range("N122").ClearContents()

# LTW row 126
$ws_LTW.Range("H126").Value = 4561.5713
$ws_LTW.Range("I126").Value = 4386.4
$ws_LTW.Range("K126").Value = 13159.2
$ws_LTW.Range("M126").Value = -10689.2

# WVR row 122
$ws_WVR.Range("H122").Value = 4312.5356
$ws_WVR.Range("I122").Value = 4388.278
$ws_WVR.Range("K122").Value = 13164.834
$ws_WVR.Range("M122").Value = -10714.834

# WVR row 126
$ws_WVR.Range("H126").Value = 1000.5263
$ws_WVR.Range("I126").Value = 1000.55554
$ws_WVR.Range("K126").Value = 3001.66662
$ws_WVR.Range("M126").Value = -531.66662
